# Update Dmanisi disability_prevalence.xlsx:
#  - New report title referencing the Unified database of targeted social assistance program
#  - Split the single "disability persons" row into two rows:
#      "family with disabilities Persons " and "disabilities Persons "
#  - New data values for 2017-2024 for both rows
#  - Minor style/border tweaks that come along with the new layout

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Make room for the extra data row: insert one blank row above the old
#    data row (row 4). This pushes the old data row to 5 and the old
#    "Source" row to 6, while rows 1-3 keep their row numbers.
# ---------------------------------------------------------------------------
$ws.Rows("4:4").Insert()

# ---------------------------------------------------------------------------
# 2. Write the new text labels first (row label order matches how the
#    strings were introduced: row 4 label, row 5 label, then the title).
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = "family with disabilities Persons "
$ws.Range("A5").Value = "disabilities Persons "
$ws.Range("A1").Value = "The number of persons with disabilities registered in the Unified database of targeted social assistance program in Dmanisi Municipality"

# ---------------------------------------------------------------------------
# 3. Row 1 - title. Widen the merge to the full table width (A1:I1) and
#    center/wrap it.
# ---------------------------------------------------------------------------
$ws.Range("A1:I1").Merge()
$ws.Range("A1:I1").Font.Name = "Arial"
$ws.Range("A1:I1").Font.Size = 11
$ws.Range("A1:I1").Font.Bold = $true
$ws.Range("A1:I1").HorizontalAlignment = -4108   # xlCenter
$ws.Range("A1:I1").VerticalAlignment = -4108     # xlCenter
$ws.Range("A1:I1").WrapText = $true
$ws.Rows(1).RowHeight = 51

# ---------------------------------------------------------------------------
# 4. Row 2 - "(End of year, persons)" caption keeps its text; it reverts to
#    the sheet's default (auto) row height instead of a custom one.
# ---------------------------------------------------------------------------
$ws.Rows(2).AutoFit()

# ---------------------------------------------------------------------------
# 5. Row 3 - year header row. A3 (empty corner cell) switches to the
#    Sylfaen font; the year cells (B3:I3) are untouched.
# ---------------------------------------------------------------------------
$ws.Range("A3").Font.Name = "Sylfaen"
$ws.Range("A3").Font.Size = 11

# ---------------------------------------------------------------------------
# 6. Row 4 (new) - "family with disabilities Persons " + first data series.
# ---------------------------------------------------------------------------
$ws.Range("A4").Font.Name = "Arial"
$ws.Range("A4").Font.Size = 10
$ws.Range("A4").Font.Bold = $false
$ws.Range("A4").Interior.Pattern = 1            # xlSolid
$ws.Range("A4").Interior.ThemeColor = 1
$ws.Range("A4").Interior.TintAndShade = 0
$ws.Range("A4").HorizontalAlignment = -4131     # xlLeft
$ws.Range("A4").VerticalAlignment = -4108       # xlCenter
$ws.Range("A4").WrapText = $true
$ws.Range("A4").Borders(8).LineStyle = 1        # xlEdgeTop
$ws.Range("A4").Borders(8).Weight = 2           # xlThin

$data4 = @(282, 268, 284, 306, 322, 335, 349, 347)
$cols = @("B", "C", "D", "E", "F", "G", "H", "I")
for ($i = 0; $i -lt 8; $i++) {
    $cell = $ws.Range($cols[$i] + "4")
    $cell.Value = $data4[$i]
    $cell.NumberFormat = "#\ ##0"
    $cell.Font.Name = "Arial"
    $cell.Font.Size = 10
    $cell.Interior.Pattern = 1
    $cell.Interior.ThemeColor = 1
    $cell.Interior.TintAndShade = 0
}
$ws.Rows(4).RowHeight = 24.75

# ---------------------------------------------------------------------------
# 7. Row 5 (was the old data row) - "disabilities Persons " + second data
#    series; only a bottom border remains under this row.
# ---------------------------------------------------------------------------
$ws.Range("A5").Font.Name = "Arial"
$ws.Range("A5").Font.Size = 10
$ws.Range("A5").Font.Bold = $false
$ws.Range("A5").Interior.Pattern = 1
$ws.Range("A5").Interior.ThemeColor = 1
$ws.Range("A5").Interior.TintAndShade = 0
$ws.Range("A5").HorizontalAlignment = -4131
$ws.Range("A5").VerticalAlignment = -4108
$ws.Range("A5").WrapText = $true
$ws.Range("A5").Borders(8).LineStyle = -4142    # xlEdgeTop off
$ws.Range("A5").Borders(9).LineStyle = 1        # xlEdgeBottom
$ws.Range("A5").Borders(9).Weight = 2

$data5 = @(320, 306, 326, 349, 366, 381, 394, 390)
for ($i = 0; $i -lt 8; $i++) {
    $cell = $ws.Range($cols[$i] + "5")
    $cell.Value = $data5[$i]
    $cell.NumberFormat = "#\ ##0"
    $cell.Font.Name = "Arial"
    $cell.Font.Size = 10
    $cell.Interior.Pattern = 1
    $cell.Interior.ThemeColor = 1
    $cell.Interior.TintAndShade = 0
    $cell.Borders(8).LineStyle = -4142
    $cell.Borders(9).LineStyle = -4142
}
# Only the last cell of the row (I5) keeps a bottom border.
$ws.Range("I5").Borders(9).LineStyle = 1
$ws.Range("I5").Borders(9).Weight = 2
$ws.Rows(5).RowHeight = 21

# ---------------------------------------------------------------------------
# 8. Row 6 - "Source" row (was row 5, shifted down automatically). Only the
#    referenced text changed in the original edit (unchanged here), so just
#    restore the row height.
# ---------------------------------------------------------------------------
$ws.Rows(6).RowHeight = 27.75

# ---------------------------------------------------------------------------
# 9. Column A width.
# ---------------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 20

Write-Output "edit complete"
